$d = $word.ActiveDocument

# --- Step 1: remove whole paragraphs that disappear entirely in the target ---
# (delete from the bottom up so earlier indices remain valid)
$toDelete = @(10, 9, 8, 7, 2)
foreach ($idx in $toDelete) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
}

# --- Step 2: helper to build a pkg:package InsertXML payload wrapping a <w:p> ---
function New-ParagraphXml([string]$innerRuns) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerRuns + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$symRun = '<w:r><w:sym w:font="Wingdings" w:char="F0E8"/></w:r>'

# --- Paragraph 1: login -> add Wingdings arrow + " Tanvi" ---
$p1 = $d.Paragraphs.Item(1)
$inner1 = '<w:r><w:t xml:space="preserve">As a logged-out user, I want to be able to login (log-in function from firebase) </w:t></w:r>' + $symRun + '<w:r><w:t xml:space="preserve"> Tanvi</w:t></w:r>'
[void]$p1.Range.InsertXML((New-ParagraphXml $inner1))

# --- Paragraph 2: browse exercise collection -> add Wingdings arrow + " José" ---
$p2 = $d.Paragraphs.Item(2)
$inner2 = '<w:r><w:t xml:space="preserve">As a logged-out user, I want to be able to browse from an exercise collection (filtering, create a database and fetch the data) </w:t></w:r>' + $symRun + '<w:r><w:t xml:space="preserve"> Jos' + [char]0x00E9 + '</w:t></w:r>'
[void]$p2.Range.InsertXML((New-ParagraphXml $inner2))

# --- Paragraph 3: create account -> add Wingdings arrow + " Tanvi" ---
$p3 = $d.Paragraphs.Item(3)
$inner3 = '<w:r><w:t xml:space="preserve">As a logged-out user, I want to have the possibility to create an account (account registration form from firebase) </w:t></w:r>' + $symRun + '<w:r><w:t xml:space="preserve"> Tanvi</w:t></w:r>'
[void]$p3.Range.InsertXML((New-ParagraphXml $inner3))

# --- Paragraph 4: save favorite exercises -> add Wingdings arrow + " Nico" ---
$p4 = $d.Paragraphs.Item(4)
$inner4 = '<w:r><w:t xml:space="preserve">As a logged-in user, I want to be able to save my favorite exercises (creation of personalized dataset using user id and display the favorite exercises chosen) </w:t></w:r>' + $symRun + '<w:r><w:t xml:space="preserve"> Nico</w:t></w:r>'
[void]$p4.Range.InsertXML((New-ParagraphXml $inner4))

# --- Paragraph 5: log workout data -> append visualize-workouts sentence + Wingdings arrow + " Nico" ---
$p5 = $d.Paragraphs.Item(5)
$inner5 = '<w:r><w:t xml:space="preserve">As a logged-in user, I want to be able to log my workout data (database using user id and saving the data into an object) </w:t></w:r><w:r><w:t>be able to visualize workouts by selecting a specific date</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>' + $symRun + '<w:r><w:t xml:space="preserve"> Nico</w:t></w:r>'
[void]$p5.Range.InsertXML((New-ParagraphXml $inner5))

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "[$i]" $p.Range.Text
}
